# ---------------------------------------------------------------------------
# Applies the "add Status column + re-style header" edit to users.xlsx
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column C (DNI) -- becomes the "Status" column.
#    Everything from the old C column onward shifts one column to the right.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).Insert()

# Give the new column its target width (raw OOXML width unit 4).
$ws.Columns.Item(3).ColumnWidth = 3.1640625

# ---------------------------------------------------------------------------
# 2. The old "Total: 35" merged banner (B3:C3) became B3:D3 after the column
#    insert (Excel auto-grows merges that straddle an inserted column).
#    Target layout: no merge, B3/C3 blank + centered, D3 untouched (no cell).
# ---------------------------------------------------------------------------
$ws.Range("B3:D3").UnMerge()
$ws.Range("B3").ClearContents()
$ws.Range("D3").Clear()

$ws.Range("B3:C3").HorizontalAlignment = -4108  #  xlCenter

# ---------------------------------------------------------------------------
# 3. New header text in the new column + row 5 shifts automatically because
#    of the column insert performed above.
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "Status"

# ---------------------------------------------------------------------------
# 4. Re-style the header row (B5:R5): navy fill, white font, centered,
#    medium outer / thin inner borders.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("B5:R5")
$headerRange.Interior.Color = 6299649  #  RGB(1,32,96)  -> FF012060
$headerRange.Font.Color = 16777215  #  RGB(255,255,255) -> white
$headerRange.Font.Name = "Abadi"
$headerRange.Font.Bold = $false
$headerRange.HorizontalAlignment = -4108  #  xlCenter
$headerRange.VerticalAlignment = -4108  #  xlCenter

$headerRange.Borders.Item(9).LineStyle = 1  #  xlEdgeBottom
$headerRange.Borders.Item(9).Weight = -4138  #  xlMedium
$headerRange.Borders.Item(8).LineStyle = 1  #  xlEdgeTop
$headerRange.Borders.Item(8).Weight = -4138
$headerRange.Borders.Item(11).LineStyle = 1  #  xlInsideVertical
$headerRange.Borders.Item(11).Weight = 2  #  xlThin

$ws.Range("B5").Borders.Item(7).LineStyle = 1  #  xlEdgeLeft
$ws.Range("B5").Borders.Item(7).Weight = -4138
$ws.Range("R5").Borders.Item(10).LineStyle = 1  #  xlEdgeRight
$ws.Range("R5").Borders.Item(10).Weight = -4138

# "Status" is displayed rotated 90 degrees in its narrow column.
$ws.Range("C5").Orientation = 90

# ---------------------------------------------------------------------------
# 5. Row heights: a thin separator row 4, and a tall header row 5.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 129.75

# ---------------------------------------------------------------------------
# 6. View state: zoom 100%, selection on H10.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("H10").Select()
